$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 199.5
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2999
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2968.4731
$ws.Range("J138").Value = 2969.137
$ws.Range("L138").Value = 8907.411
$ws.Range("N138").Value = -19187.411

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 1000150
$ws.Range("I8").Value = 1000150
$ws.Range("K8").Value = 1000150
$ws.Range("M8").Value = -1000006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2000000
$ws.Range("I11").Value = 2000000
$ws.Range("K11").Value = 2000000
$ws.Range("M11").Value = -1999856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 334010000
$ws.Range("I13").Value = 501000000
$ws.Range("J13").Value = 30000
$ws.Range("K13").Value = 501000000
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = -500999856
$ws.Range("N13").Value = -30288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5021.6763
$ws.Range("I32").Value = 4741.7812
$ws.Range("K32").Value = 4741.7812
$ws.Range("M32").Value = -4454.7812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5224.1816

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3014.2144
$ws.Range("I74").Value = 2271.3333
$ws.Range("K74").Value = 2271.3333
$ws.Range("M74").Value = -1397.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 9288
$ws.Range("J76").Value = 9288
$ws.Range("L76").Value = 9288
$ws.Range("N76").Value = -9964

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3014.2144
$ws.Range("I77").Value = 2271.3333
$ws.Range("K77").Value = 11356.6665
$ws.Range("M77").Value = -6988.666499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 9288
$ws.Range("J79").Value = 9288
$ws.Range("L79").Value = 9288
$ws.Range("N79").Value = -11628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1245.5
$ws.Range("I102").Value = 1327.3334
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 1327.3334
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 294.6666
$ws.Range("N102").Value = -4244

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 39999
$ws.Range("J104").Value = 39999
$ws.Range("L104").Value = 39999
$ws.Range("N104").Value = -46987

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 132391.58
$ws.Range("J113").Value = 132391.58
$ws.Range("L113").Value = 132391.58
$ws.Range("N113").Value = -141069.58

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 923.1053000000001
$ws.Range("I122").Value = 818.8333
$ws.Range("K122").Value = 2456.4999
$ws.Range("M122").Value = -6.499899999999798

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1018.03845
$ws.Range("I132").Value = 869
$ws.Range("K132").Value = 2607
$ws.Range("M132").Value = -77

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5224.1816

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 34330
$ws.Range("J6").Value = 34330
$ws.Range("L6").Value = 34330
$ws.Range("N6").Value = -34556

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1949.75
$ws.Range("I20").Value = 1949.75
$ws.Range("K20").Value = 1949.75
$ws.Range("M20").Value = -1702.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5303.5
$ws.Range("I107").Value = 4761.6
$ws.Range("K107").Value = 4761.6
$ws.Range("M107").Value = -2841.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 439.33334
$ws.Range("I134").Value = 439.33334
$ws.Range("K134").Value = 1318.00002
$ws.Range("M134").Value = 1216.99998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1956.625
$ws.Range("I31").Value = 1700.9
$ws.Range("J31").Value = 3235.25
$ws.Range("K31").Value = 1700.9
$ws.Range("L31").Value = 3235.25
$ws.Range("M31").Value = -1405.9
$ws.Range("N31").Value = -3825.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1956.625
$ws.Range("I34").Value = 1700.9
$ws.Range("J34").Value = 3235.25
$ws.Range("K34").Value = 1700.9
$ws.Range("L34").Value = 3235.25
$ws.Range("M34").Value = -1498.9
$ws.Range("N34").Value = -3639.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 2500
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 10708.25
$ws.Range("J88").Value = 10708.25
$ws.Range("L88").Value = 10708.25
$ws.Range("N88").Value = -11520.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 10708.25
$ws.Range("J91").Value = 10708.25
$ws.Range("L91").Value = 10708.25
$ws.Range("N91").Value = -13516.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1175.2
$ws.Range("I122").Value = 893.38464
$ws.Range("J122").Value = 3007
$ws.Range("K122").Value = 2680.15392
$ws.Range("L122").Value = 9021
$ws.Range("M122").Value = -230.1539199999997
$ws.Range("N122").Value = -13921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 22000090
$ws.Range("I4").Value = 22000090
$ws.Range("K4").Value = 66000270
$ws.Range("M4").Value = -66000158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 15640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1053.6316
$ws.Range("I113").Value = 1829.6666
$ws.Range("J113").Value = 908.125
$ws.Range("K113").Value = 5488.9998
$ws.Range("L113").Value = 2724.375
$ws.Range("M113").Value = -3318.9998
$ws.Range("N113").Value = -7064.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2543.2666
$ws.Range("I131").Value = 1699.8572
$ws.Range("J131").Value = 3281.25
$ws.Range("K131").Value = 5099.571599999999
$ws.Range("L131").Value = 9843.75
$ws.Range("M131").Value = -59.57159999999931
$ws.Range("N131").Value = -19923.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3930.2
$ws.Range("I132").Value = 2607.4285
$ws.Range("J132").Value = 5613.727
$ws.Range("K132").Value = 23466.8565
$ws.Range("L132").Value = 50523.543
$ws.Range("M132").Value = -20936.8565
$ws.Range("N132").Value = -55583.543

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 16000
$ws.Range("I133").Value = 4000
$ws.Range("K133").Value = 12000
$ws.Range("M133").Value = -6940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3153.8572
$ws.Range("J137").Value = 4026.3333
$ws.Range("L137").Value = 12078.9999
$ws.Range("N137").Value = -22278.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 551.5
$ws.Range("I97").Value = 489.25
$ws.Range("K97").Value = 489.25
$ws.Range("M97").Value = 6.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1557.3334
$ws.Range("I107").Value = 335.46155
$ws.Range("K107").Value = 335.46155
$ws.Range("M107").Value = 1584.53845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1585.5264
$ws.Range("I132").Value = 1789.1875
$ws.Range("J132").Value = 499.33334
$ws.Range("K132").Value = 5367.5625
$ws.Range("L132").Value = 1498.00002
$ws.Range("M132").Value = -2837.5625
$ws.Range("N132").Value = -6558.000019999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 951.8570999999999
$ws.Range("I22").Value = 773.53845
$ws.Range("K22").Value = 773.53845
$ws.Range("M22").Value = -478.53845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 66666.336
$ws.Range("I23").Value = 66666.336
$ws.Range("K23").Value = 66666.336
$ws.Range("M23").Value = -66436.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 951.8570999999999
$ws.Range("I27").Value = 773.53845
$ws.Range("K27").Value = 773.53845
$ws.Range("M27").Value = -666.53845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 16667800
$ws.Range("I100").Value = 50000376
$ws.Range("J100").Value = 1512.5
$ws.Range("K100").Value = 100000752
$ws.Range("L100").Value = 3025
$ws.Range("M100").Value = -100000211
$ws.Range("N100").Value = -4107

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 72013.5
$ws.Range("J117").Value = 72013.5
$ws.Range("L117").Value = 72013.5
$ws.Range("N117").Value = -81191.5
